$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the "period" values (column G) to the new short labels ---
$ws.Range("G2").Value = "evening"
$ws.Range("G3").Value = "night"

# --- 2. Add new "song" / "artist" columns derived from the existing
#        "song-artist" column (H), splitting on " - " ---
$ws.Range("J1").Value = "song"
$ws.Range("K1").Value = "artist"

$songArtist2 = $ws.Range("H2").Value
$parts2 = $songArtist2 -split " - "
$ws.Range("J2").Value = $parts2[0] + " "
$ws.Range("K2").Value = " " + $parts2[1]

$songArtist3 = $ws.Range("H3").Value
$parts3 = $songArtist3 -split " - "
$ws.Range("J3").Value = $parts3[0] + " "
$ws.Range("K3").Value = " " + $parts3[1]

# --- 3. Autofit the new / affected columns, mirroring the bestFit widths
#        that Excel itself computes for new data ---
$ws.Columns.Item(8).AutoFit() | Out-Null
$ws.Columns.Item(11).AutoFit() | Out-Null

# --- 4. Update the active selection like the author left it ---
$ws.Range("J3").Select()
